$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.608902839520766, 1.2752471729130956, -0.008016812747032125, 2.8546286833102936)"
$ws.Range("C2").Value = "EXN(2.431317568620919, 2.2002328460886167, 4.037413153604303)"
$ws.Range("D2").Value = "NIG(0.9628236769974599, 0.6470768124984324, 1.322808158597693, 3.22329903066087)"
$ws.Range("E2").Value = "NIG(0.8487520760179528, 0.5550835306140275, 5.121829312611803, 5.185113448159607)"
